# Apply the latest crypto price/volume snapshot scraped by the GitHub Action.
# Numeric-looking text values are written with a leading apostrophe so Excel
# keeps them as text (matching the original sheet, which stores these as strings)
# instead of silently re-interpreting them as numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.239.90'
$ws.Range('D3').Value = '1.862.55'
$ws.Range('E3').Value = '  -0.41%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '''236.77'
$ws.Range('E5').Value = '  +0.43%  '
$ws.Range('D6').Value = '''1.001'
$ws.Range('E6').Value = '  +0.02%  '
$ws.Range('D7').Value = '''0.4700'
$ws.Range('E7').Value = '  +0.69%  '
$ws.Range('D8').Value = '''0.2900'
$ws.Range('E8').Value = '  +2.19%  '
$ws.Range('D9').Value = '''0.06539'
$ws.Range('D10').Value = '''21.54'
$ws.Range('E10').Value = '  +1.93%  '
$ws.Range('D11').Value = '''0.07943'
$ws.Range('E11').Value = '  -0.03%  '
$ws.Range('D12').Value = '''98.13'
$ws.Range('E12').Value = '  +0.75%  '
$ws.Range('D13').Value = '1.868.81'
$ws.Range('E13').Value = '  -0.10%  '
$ws.Range('D14').Value = '''5.160'
$ws.Range('E14').Value = '  +0.10%  '
$ws.Range('D15').Value = '''0.6818'
$ws.Range('E15').Value = '  +1.08%  '
$ws.Range('D16').Value = '''267.82'
$ws.Range('E16').Value = '  -5.21%  '
$ws.Range('D17').Value = '30.233.76'
$ws.Range('E17').Value = '  -0.47%  '
$ws.Range('D18').Value = '''13.74'
$ws.Range('E18').Value = '  +8.26%  '
$ws.Range('E19').Value = '  -0.04%  '
$ws.Range('D20').Value = '''0.000007376'
$ws.Range('E20').Value = '  +1.13%  '
$ws.Range('E21').Value = '  -0.28%  '
$ws.Range('E22').Value = '  -4.12%  '
$ws.Range('D23').Value = '''1.001'
$ws.Range('E23').Value = '  +0.00%  '
$ws.Range('D24').Value = '''6.188'
$ws.Range('E24').Value = '  -0.26%  '
$ws.Range('D25').Value = '''166.60'
$ws.Range('E25').Value = '  +0.99%  '
$ws.Range('D26').Value = '''9.214'
$ws.Range('D27').Value = '''18.89'
$ws.Range('E27').Value = '  -1.11%  '
$ws.Range('E28').Value = '  +1.18%  '
$ws.Range('D29').Value = '''1.392'
$ws.Range('E29').Value = '  +2.73%  '
$ws.Range('D30').Value = '''0.09826'
$ws.Range('E30').Value = '  +1.28%  '
$ws.Range('D31').Value = '''4.370'
$ws.Range('E31').Value = '  -1.62%  '
$ws.Range('D32').Value = '''1.472'
$ws.Range('E32').Value = '  -0.33%  '
$ws.Range('D33').Value = '''4.045'
$ws.Range('E33').Value = '  -1.71%  '
$ws.Range('D34').Value = '''0.04703'
$ws.Range('E34').Value = '  -0.10%  '
$ws.Range('D35').Value = '''1.129'
$ws.Range('E35').Value = '  +0.89%  '
$ws.Range('D36').Value = '''0.7039'
$ws.Range('E36').Value = '  -0.13%  '
$ws.Range('D37').Value = '''2.709'
$ws.Range('E37').Value = '  -0.38%  '
$ws.Range('E38').Value = '  +0.56%  '
$ws.Range('E39').Value = '  +2.56%  '
$ws.Range('D40').Value = '''6.282'
$ws.Range('E40').Value = '  -0.78%  '
$ws.Range('D41').Value = '''74.35'
$ws.Range('E41').Value = '  +1.06%  '
$ws.Range('D42').Value = '''1.935'
$ws.Range('E42').Value = '  -0.67%  '
$ws.Range('D43').Value = '''0.8456'
$ws.Range('D44').Value = '''0.4163'
$ws.Range('E44').Value = '  -0.91%  '
$ws.Range('E45').Value = '  -0.08%  '
$ws.Range('D46').Value = '''102.95'
$ws.Range('E46').Value = '  -0.89%  '
$ws.Range('D47').Value = '''954.07'
$ws.Range('E47').Value = '  +1.54%  '
$ws.Range('D48').Value = '''7.148'
$ws.Range('E48').Value = '  -0.92%  '
$ws.Range('D49').Value = '''9.218'
$ws.Range('E49').Value = '  -0.15%  '
$ws.Range('D50').Value = '''34.13'
$ws.Range('E50').Value = '  -0.17%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').Value = '''0.05653'
$ws.Range('E51').Value = '  +0.48%  '
